# Update vm_pu.xlsx results sheet with new values (case with 380 kV done)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that receive new values, in order: B, C, D, E, F, I, J, K, L, M, N
# (column G stays at 1, column H is empty in every data row, both unchanged)
$cols = @(2, 3, 4, 5, 6, 9, 10, 11, 12, 13, 14)

# Each inner array: row number followed by the 11 new values for columns listed above
$data = @(
    @(2, 1.02, 1.036366119186613, 1.03675451196062, 1.049606669417961, 1.056448385976435, 1.033266918840126, 1.041474828900946, 1.039547348927932, 1.052363258609685, 1.059186068828082, 1.042953842616818),
    @(3, 1.02, 1.037513624083608, 1.037579487872779, 1.050777947643711, 1.05776394337133, 1.033480814666864, 1.042265270926398, 1.040181951009288, 1.053345890141792, 1.060313989994319, 1.043745407160612),
    @(4, 1.02, 1.038255734688324, 1.03811248070693, 1.051535847457616, 1.058615471772264, 1.033617091008566, 1.042775812769769, 1.04059111487002, 1.053981148218219, 1.061043558205735, 1.04425667403196),
    @(5, 1.02, 1.038567623322632, 1.038336354339919, 1.051854470842395, 1.058973521760297, 1.033673871953995, 1.042990223389462, 1.040762776678862, 1.054248075256152, 1.06135020486955, 1.044471389139326),
    @(6, 1.02, 1.03861998528982, 1.038373932186733, 1.051907969281386, 1.059033643932155, 1.033683375838186, 1.043026210928163, 1.040791578896891, 1.054292885567655, 1.061401688490965, 1.044507427784459),
    @(7, 1.02, 1.038259902529592, 1.038115472889513, 1.051540104908309, 1.058620255786642, 1.033617851720912, 1.042778678603095, 1.040593410001078, 1.05398471544132, 1.061047655882531, 1.044259543935098),
    @(8, 1.02, 1.03675400795476, 1.03703348668892, 1.05000250819363, 1.056892928860897, 1.033339646643573, 1.04174215488349, 1.03976211943519, 1.052695462245807, 1.059567312246512, 1.04322154823292),
    @(9, 1.02, 1.034097282994934, 1.035120586224591, 1.047293044777958, 1.053851187693952, 1.032833112822083, 1.039908533325197, 1.038286035434325, 1.050419217990098, 1.056956611499878, 1.041385322721695),
    @(10, 1.02, 1.032323918061507, 1.033841065450837, 1.04548663998973, 1.051824616758407, 1.032484463323707, 1.038681270974726, 1.037294395419505, 1.048898682844533, 1.055214612546413, 1.040156317517917),
    @(11, 1.02, 1.031555483021681, 1.03328600507769, 1.044704404938441, 1.050947362188706, 1.032330893789302, 1.038148690653748, 1.036863199110112, 1.048239536948772, 1.054459923725731, 1.039622980871796),
    @(12, 1.02, 1.031269966350654, 1.033079677238012, 1.044413838761693, 1.050621548182855, 1.032273460139957, 1.03795069001994, 1.036702761052387, 1.047994587474765, 1.054179537937757, 1.039424699054377),
    @(13, 1.02, 1.031331214556985, 1.033123942224298, 1.044476166639893, 1.05069143473712, 1.032285797552943, 1.037993169837916, 1.036737187904738, 1.048047135122161, 1.054239684421053, 1.039467239198567),
    @(14, 1.02, 1.031531883884682, 1.033268953092794, 1.044680386846535, 1.050920429551512, 1.032326154284967, 1.03813232746875, 1.036849942814067, 1.048219291669887, 1.054436748199477, 1.039606594449198),
    @(15, 1.02, 1.031655511463559, 1.03335827869116, 1.044806212377476, 1.051061525702511, 1.032350967571643, 1.038218043642982, 1.036919378699342, 1.048325347935922, 1.054558157589399, 1.039692432350229),
    @(16, 1.02, 1.032374904664133, 1.033877881455291, 1.045538553153522, 1.051882842624554, 1.032494600392853, 1.038716591876735, 1.037322974314394, 1.048942412377281, 1.055264690341931, 1.040191688579658),
    @(17, 1.02, 1.032826010590239, 1.034203541230458, 1.045997916588239, 1.052398102121705, 1.032584000723579, 1.03902900458596, 1.037575654219314, 1.04932928028669, 1.055707773445878, 1.040504544950764),
    @(18, 1.02, 1.033089079509506, 1.034393394744269, 1.04626585081602, 1.052698669895432, 1.032635895380184, 1.039211117053682, 1.037722863574654, 1.049554862097837, 1.055966178650981, 1.040686916039079),
    @(19, 1.02, 1.033178770195906, 1.034458113225228, 1.046357208721163, 1.052801160154326, 1.03265354754053, 1.039273193638795, 1.037773028572739, 1.049631767551845, 1.056054281815671, 1.040749080780064),
    @(20, 1.02, 1.032777616745862, 1.034168611192355, 1.045948631736776, 1.052342817019371, 1.032574434881883, 1.038995497311775, 1.037548562116875, 1.049287780447909, 1.055660238732218, 1.040470990092407),
    @(21, 1.02, 1.031472794151421, 1.033226255253722, 1.044620249355255, 1.050852995272011, 1.032314281030525, 1.03809135389904, 1.036816746832763, 1.048168598998633, 1.054378719520911, 1.039565562692319),
    @(22, 1.02, 1.030651903616047, 1.032632869114879, 1.043784988233613, 1.049916500470267, 1.032148449184191, 1.037521860772578, 1.036355048289723, 1.047464268681232, 1.053572625058251, 1.038995260820296),
    @(23, 1.02, 1.031087121014747, 1.032947518839871, 1.044227781704523, 1.050412934438926, 1.032236574327925, 1.037823857146814, 1.036599953153689, 1.047837710256895, 1.053999985031557, 1.039297686064023),
    @(24, 1.02, 1.032799484005875, 1.034184394876552, 1.045970901453687, 1.052367797897586, 1.032578758049534, 1.039010638154146, 1.037560804410867, 1.049306532659769, 1.055681717739066, 1.040486152436511),
    @(25, 1.02, 1.034784493144976, 1.035615865355015, 1.047993516906357, 1.054637321268881, 1.032965995712148, 1.040383419142842, 1.038668973433936, 1.051008212138918, 1.057631803657777, 1.041860882931671)
)

foreach ($rowData in $data) {
    $r = $rowData[0]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $c = $cols[$i]
        $val = $rowData[$i + 1]
        $ws.Cells.Item($r, $c).Value = $val
    }
}
